$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A17:C17").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A18").Value = "Progress bar or output for PDF creator"
$ws.Range("B18").Value = "When clicking the ""Create PDF"" button in the document viewer, the user is not given any indication of what is happening. It would be nice if there were a progress bar or output window showing what is happening."
$ws.Range("C18").Value = "OPEN"

$ws.Rows.Item(18).RowHeight = 30

$ws.Range("C18").Select()
